# edit.ps1 -- applies the "fixed reset widget entries upon row insertion" change:
#   * "2022" sheet gains 20 new data rows (7-26), some rows only carry a
#     partial widget snapshot (a subset of columns), mirroring what the app
#     had written before the reset-bug fix.
#   * "2023" sheet gains one new data row (3).
#   * "2024" sheet: row 2 is corrected (OBR#/Category/Brand/Price/Notes) and
#     a new row 3 is appended; the sheet becomes the active tab with G2
#     selected.
#
# Helper: writes a cell as TEXT, even when the text looks like a number
# (e.g. "888888") or is empty - matches how the source workbook stores its
# data (t="inlineStr"/shared-string, never auto-converted to a number) and
# avoids leaving a stray NumberFormat/quotePrefix style behind.
function Set-TextCell {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "2022" -- append rows 7-26
# ---------------------------------------------------------------------------
$ws2022 = $wb.Worksheets.Item("2022")

$rows2022 = @(
    @{ Row = 7;  A = "test month";       B = "888888";   C = "Table";           D = "Executive";  E = "zooey";       F = "123123";     G = "notes dito" }
    @{ Row = 8;  A = "3";                B = "321321";   C = "Others";          D = "Sound System"; E = "jbl";       F = "1000000000"; G = "notes?" }
    @{ Row = 9;  A = "test ulit";        B = "asdf";     C = "Kitchen-related"; D = "Microwave Oven"; E = "asd";     F = "a";          G = "a" }
    @{ Row = 10; A = "111";              B = "11111";    C = "Chair";           D = "Executive";                     F = "asf";        G = "a" }
    @{ Row = 11;                                         C = "PY_VAR0";         D = "Executive" }
    @{ Row = 12;                                         C = "Kitchen-related"; D = "Microwave Oven" }
    @{ Row = 13; A = "month focus line";                 C = "Others";          D = "AirCon" }
    @{ Row = 14; A = "month";            B = "focus";    C = "Kitchen-related"; D = "Rice Cooker"; E = "asdf";       F = "asdf" }
    @{ Row = 15; A = "trying-out";       B = "auto set"; C = "Computer";        D = "AVR / UPS";   E = "instead of"; F = "py_var0" }
    @{ Row = 16;                                         C = "Select Category" }
    @{ Row = 17;                                         C = "Select Category" }
    @{ Row = 18;                                         C = "Select Category" }
    @{ Row = 19;                                         C = "Select Category" }
    @{ Row = 20;                                         C = "Select Category" }
    @{ Row = 21;                                         C = "Select Category" }
    @{ Row = 22;                                         C = "Select Category" }
    @{ Row = 23;                                         C = "Select Category" }
    @{ Row = 24;                                         C = "Select Category" }
    @{ Row = 25;                                         C = "Select Category" }
    @{ Row = 26; A = "";                 B = "";         C = "Select Category"; D = "";            E = "";           F = "";            G = "" }
)

foreach ($r in $rows2022) {
    $rowNum = $r.Row
    foreach ($col in @("A", "B", "C", "D", "E", "F", "G")) {
        if ($r.ContainsKey($col)) {
            Set-TextCell $ws2022.Range("$col$rowNum") $r[$col]
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "2023" -- append row 3
# ---------------------------------------------------------------------------
$ws2023 = $wb.Worksheets.Item("2023")
Set-TextCell $ws2023.Range("A3") "test month"
Set-TextCell $ws2023.Range("B3") "888888"
Set-TextCell $ws2023.Range("C3") "Table"
Set-TextCell $ws2023.Range("D3") "Executive"
Set-TextCell $ws2023.Range("E3") "zooey"
Set-TextCell $ws2023.Range("F3") "123123"
Set-TextCell $ws2023.Range("G3") "notes dito"

# ---------------------------------------------------------------------------
# Sheet "2024" -- fix row 2, append row 3, make it the active sheet with G2
# selected
# ---------------------------------------------------------------------------
$ws2024 = $wb.Worksheets.Item("2024")

Set-TextCell $ws2024.Range("C2") "wwqr"
Set-TextCell $ws2024.Range("E2") "brando"
Set-TextCell $ws2024.Range("G2") "departments pa na dropdown"
$ws2024.Range("B2").Value = 2024
$ws2024.Range("B2").ClearFormats()
$ws2024.Range("F2").Value = 12378141
$ws2024.Range("F2").ClearFormats()

Set-TextCell $ws2024.Range("A3") "di nagki clear all"
Set-TextCell $ws2024.Range("B3") "hahahaha"
Set-TextCell $ws2024.Range("C3") "Table"
Set-TextCell $ws2024.Range("D3") "Conference"
Set-TextCell $ws2024.Range("E3") "pang 2024 na page"
Set-TextCell $ws2024.Range("F3") "5555555"
Set-TextCell $ws2024.Range("G3") "notes dito"

$ws2024.Activate() | Out-Null
$ws2024.Range("G2").Select() | Out-Null
